$wb = $excel.ActiveWorkbook

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 37039228
$ws.Range("I132").Value = 41668520
$ws.Range("J132").Value = 4900
$ws.Range("K132").Value = 125005560
$ws.Range("L132").Value = 14700
$ws.Range("M132").Value = -125003030
$ws.Range("N132").Value = -19760

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 570127.4
$ws.Range("I32").Value = 4234.614
$ws.Range("J32").Value = 3336714
$ws.Range("K32").Value = 4234.614
$ws.Range("L32").Value = 3336714
$ws.Range("M32").Value = -3947.614
$ws.Range("N32").Value = -3337288

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3233.0952
$ws.Range("I74").Value = 512.1875
$ws.Range("J74").Value = 11940
$ws.Range("K74").Value = 512.1875
$ws.Range("L74").Value = 11940
$ws.Range("M74").Value = 361.8125
$ws.Range("N74").Value = -13688

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3233.0952
$ws.Range("I77").Value = 512.1875
$ws.Range("J77").Value = 11940
$ws.Range("K77").Value = 2560.9375
$ws.Range("L77").Value = 59700
$ws.Range("M77").Value = 1807.0625
$ws.Range("N77").Value = -68436

# BSM row 62
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 30192.334
$ws.Range("J62").Value = 30192.334
$ws.Range("L62").Value = 30192.334
$ws.Range("N62").Value = -31564.334

# BSM row 65
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 30192.334
$ws.Range("J65").Value = 30192.334
$ws.Range("L65").Value = 90577.00199999999
$ws.Range("N65").Value = -97441.00199999999

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1908.3226
$ws.Range("I134").Value = 1735.75
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 5207.25
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -2672.25
$ws.Range("N134").Value = -12570

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 356.41177
$ws.Range("I4").Value = 84.85714
$ws.Range("J4").Value = 546.5
$ws.Range("K4").Value = 254.57142
$ws.Range("L4").Value = 1639.5
$ws.Range("M4").Value = -142.57142
$ws.Range("N4").Value = -1863.5

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 850
$ws.Range("I5").Value = 700
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 2100
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1988
$ws.Range("N5").Value = -3224

# CUL row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 117.166664
$ws.Range("I6").Value = 34.333332
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 102.999996
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = 10.000004
$ws.Range("N6").Value = -826

# CUL row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 449.4
$ws.Range("I7").Value = 373.5
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 1120.5
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -1008.5
$ws.Range("N7").Value = -1724

# CUL row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 62557084
$ws.Range("I9").Value = 250000290
$ws.Range("J9").Value = 76016.836
$ws.Range("K9").Value = 750000870
$ws.Range("L9").Value = 228050.508
$ws.Range("M9").Value = -750000646
$ws.Range("N9").Value = -228498.508

# CUL row 10
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 255.93333
$ws.Range("I10").Value = 83.90000000000001
$ws.Range("K10").Value = 251.7
$ws.Range("M10").Value = -112.7

# CUL row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 956.4167
$ws.Range("I11").Value = 1561.4286
$ws.Range("J11").Value = 109.4
$ws.Range("K11").Value = 4684.2858
$ws.Range("L11").Value = 328.2
$ws.Range("M11").Value = -4544.2858
$ws.Range("N11").Value = -608.2

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 92.166664
$ws.Range("I12").Value = 167.4
$ws.Range("K12").Value = 502.2
$ws.Range("M12").Value = -329.2

# CUL row 13
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 195
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 195
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 585
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -921

# CUL row 15
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 140.8
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 140.8
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 422.4
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -702.4000000000001

# CUL row 16
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 700
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2100
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2446

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 300
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 600
$ws.Range("M17").Value = -1331
$ws.Range("N17").Value = -938

# CUL row 19
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 725
$ws.Range("I19").Value = 900
$ws.Range("J19").Value = 666.6667
$ws.Range("K19").Value = 2700
$ws.Range("L19").Value = 2000.0001
$ws.Range("M19").Value = -2526
$ws.Range("N19").Value = -2348.0001

# CUL row 20
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 2480
$ws.Range("J20").Value = 2480
$ws.Range("L20").Value = 7440
$ws.Range("N20").Value = -7894

# CUL row 21
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 992.44446
$ws.Range("I21").Value = 500
$ws.Range("J21").Value = 1238.6666
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 3715.9998
$ws.Range("M21").Value = -1327
$ws.Range("N21").Value = -4061.9998

# CUL row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 737.41815
$ws.Range("J22").Value = 737.41815
$ws.Range("L22").Value = 2212.25445
$ws.Range("N22").Value = -2550.25445

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 53.5
$ws.Range("I23").Value = 42.25
$ws.Range("K23").Value = 126.75
$ws.Range("M23").Value = 108.25

# CUL row 24
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 370
$ws.Range("J24").Value = 433.33334
$ws.Range("L24").Value = 1300.00002
$ws.Range("N24").Value = -1760.00002

# CUL row 25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 755
$ws.Range("I25").Value = 1500
$ws.Range("J25").Value = 506.66666
$ws.Range("K25").Value = 4500
$ws.Range("L25").Value = 1519.99998
$ws.Range("M25").Value = -4331
$ws.Range("N25").Value = -1857.99998

# CUL row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 5952645
$ws.Range("I26").Value = 149.33333
$ws.Range("J26").Value = 10417016
$ws.Range("K26").Value = 447.99999
$ws.Range("L26").Value = 31251048
$ws.Range("M26").Value = -159.99999
$ws.Range("N26").Value = -31251624

# CUL row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 737.41815
$ws.Range("J27").Value = 737.41815
$ws.Range("L27").Value = 2212.25445
$ws.Range("N27").Value = -2416.25445

# CUL row 29
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 75.75
$ws.Range("I29").Value = 13
$ws.Range("J29").Value = 96.666664
$ws.Range("K29").Value = 39
$ws.Range("L29").Value = 289.999992
$ws.Range("M29").Value = 238
$ws.Range("N29").Value = -843.999992

# CUL row 30
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 755
$ws.Range("I30").Value = 1500
$ws.Range("J30").Value = 506.66666
$ws.Range("K30").Value = 4500
$ws.Range("L30").Value = 1519.99998
$ws.Range("M30").Value = -4398
$ws.Range("N30").Value = -1723.99998

# CUL row 31
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1500
$ws.Range("I31").Value = 500
$ws.Range("J31").Value = 1900
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 5700
$ws.Range("M31").Value = -1212
$ws.Range("N31").Value = -6276

# CUL row 32
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1001833.5
$ws.Range("I32").Value = 1000001
$ws.Range("J32").Value = 1002000.06
$ws.Range("K32").Value = 3000003
$ws.Range("L32").Value = 3006000.18
$ws.Range("M32").Value = -2999720
$ws.Range("N32").Value = -3006566.18

# CUL row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 228.41176
$ws.Range("I33").Value = 129.125
$ws.Range("K33").Value = 774.75
$ws.Range("M33").Value = -491.75

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 850
$ws.Range("I135").Value = 700
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 6300
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -3765
$ws.Range("N135").Value = -14070

# GSM row 42
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 16000
$ws.Range("J42").Value = 16000
$ws.Range("L42").Value = 16000
$ws.Range("N42").Value = -16970

# GSM row 115
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H115").Value = 16000
$ws.Range("J115").Value = 16000
$ws.Range("L115").Value = 16000
$ws.Range("N115").Value = -18350
